# Applies the "Sample Vol (mL)" column addition described by the commit:
#  - Dilution sheet (Sheet2): insert a new column E ("Sample Vol (mL)") that
#    holds a correction formula (=D*0.9998395) derived from the existing
#    Sample Wt (g) column D. Everything right of D shifts over by one
#    column (old E..J -> F..K).
#  - Sheet1: mirror the "Sample Wt (g)" / "Total vol (mL)" headers from the
#    Dilution sheet into new columns E/F (header row only), matching the
#    highlighted-header style used on the Dilution sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Dilution sheet: insert new column E with the sample-volume formula ---
$ws2.Columns.Item(5).Insert()

$ws2.Range("E1").Value = "Sample Vol (mL)"

$lastRow = 14
for ($r = 2; $r -le $lastRow; $r++) {
    $ws2.Cells.Item($r, 5).Formula = "=D" + $r + "*0.9998395"
}

# match the column width of the donor column (D) as closely as possible
$ws2.Columns.Item(5).ColumnWidth = $ws2.Columns.Item(4).ColumnWidth

# --- Sheet1: add the two mirrored headers in columns E/F ---
$ws1.Range("E1").Value = "Sample Wt (g)"
$ws1.Range("F1").Value = "Total vol (mL)"

# copy the highlighted header format (fill + bold) from the Dilution sheet
# header row onto the new Sheet1 header cells
$ws2.Range("A1").Copy()
$ws1.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Columns.Item(5).ColumnWidth = 14.33
$ws1.Columns.Item(6).ColumnWidth = 13.5

# --- restore selections: Sheet1 becomes the active tab/selection, while
#     the Dilution sheet keeps its own last-used selection ---
$ws2.Activate() | Out-Null
$ws2.Range("E18").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("H13").Select() | Out-Null
